$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that need their partner-match counts halved (corrected data), and the
# column range that holds the numeric match counts (B..AK).
$rows = @(20, 39, 72, 90)
$lastCol = 37  # column AK

foreach ($r in $rows) {
    for ($c = 2; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -ne $null -and $val -ne 0) {
            $cell.Value2 = $val / 2
        }
    }
}

# Fix mojibake-encoded organization name text in A28
$ws.Cells.Item(28, 1).Value2 = "Consejo Nacional de Investigaciones CientÃ­ficas y TÃ©cnicas"
